$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the two target paragraphs by their distinctive text so the
# script does not depend on brittle paragraph indices.
# ------------------------------------------------------------------
$paraValidar = $null
$paraMostrar = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($paraValidar -eq $null -and $t -like "*Validar que no deje ingresar dos columnas con el mismo nombre*") {
        $paraValidar = $p
    }
    if ($paraMostrar -eq $null -and $t -like "*.Mostrar Nombre de la Columna*") {
        $paraMostrar = $p
    }
}

if ($paraValidar -eq $null) {
    throw "Could not locate the 'Validar que no deje ingresar...' paragraph"
}
if ($paraMostrar -eq $null) {
    throw "Could not locate the '.Mostrar Nombre de la Columna...' paragraph"
}

# ------------------------------------------------------------------
# Paragraph 1: "Validar que no deje ingresar dos columnas con el
# mismo nombre" — highlight the paragraph mark and the run yellow.
# ------------------------------------------------------------------
$xmlValidar = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2251881A" w14:textId="76708F1C" w:rsidR="00E758E6" w:rsidRDefault="00E758E6" w:rsidP="00144ED4"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Validar que no deje ingresar dos columnas con el mismo nombre</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$paraValidar.Range.InsertXML($xmlValidar)

# ------------------------------------------------------------------
# Paragraph 2: ".Mostrar Nombre de la Columna, Nombre del E/S y Tipo"
# — split ".Mostrar" into "." (unhighlighted) + "Mostrar"
# (highlighted), highlight the trailing run and the paragraph mark.
# ------------------------------------------------------------------
$xmlMostrar = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="375BAC80" w14:textId="0DF09DF1" w:rsidR="00E758E6" w:rsidRDefault="00763DDD" w:rsidP="006553EE"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="259" w:lineRule="auto"/><w:jc w:val="left"/><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>.</w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Mostrar</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> Nombre de la Columna, Nombre del E/S y Tipo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$paraMostrar.Range.InsertXML($xmlMostrar)

Write-Output "done"
